$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.935.06'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.557.29'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '207.70'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.488'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '21.99'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0854'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').Value = '1.777.84'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '1.555.76'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.74'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.519'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '26.921.82'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.70'
$ws.Range('D17').Style = "Normal"
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '215.52'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.35'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.22'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.01'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0471'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('E31').Value = '  +2.26%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.14'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('D34').Value = '1.416.76'
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.05'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +7.70%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.33'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +1.60%  '
$ws.Range('E40').Value = '  +2.13%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.808'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.97%  '
$ws.Range('E44').Value = '  +1.66%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '64.29'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').Value = '1.692.07'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '86.84'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0519'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0957'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.36%  '
